$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.031.91'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.35%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.315.58'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.43%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.34%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '559.70'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.15%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '186.14'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.63%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.585'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.308.17'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.48%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.186'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.54%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.585'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.13%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.65'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.73%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000270'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.99%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.64'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.41%  '

$ws.Range("B15").Value = 'BitcoinCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '630.39'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.33%  '

$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C16").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.833.00'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.90%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.15'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.51%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '66.048.39'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.26%  '

$ws.Range("E19").Value = '  -1.57%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.302.30'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.85%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.31'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.90%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.908'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.66%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '18.26'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.33%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '102.34'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +7.48%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.96'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.30%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.95'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.78%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.05'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.68%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.73'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.11%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.56'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.87%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.68'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.89%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '30.33'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.69%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.07'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.46%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.39'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.16%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.12'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.24%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '553.92'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.45%  '

$ws.Range("E36").Value = '  -0.62%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.830.48'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.61%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '57.61'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.79%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0₃0739'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.76%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '33.84'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.19%  '

$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.26'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.60%  '

$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.129'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.38%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.70'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.44%  '

$ws.Range("B45").Value = 'CoreDAO'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.22'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -15.48%  '

$ws.Range("B46").Value = 'TheGraph'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.335'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.67%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0419'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.35%  '

$ws.Range("E48").Value = '  +2.62%  '

$ws.Range("B49").Value = 'ThetaToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.61'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.07%  '

$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.129'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.18%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.998'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.14%  '
